# Update the cached "Update automatically" date/time footer field
# (placeholder type dt / ppPlaceholderDate) from 3/7/19 to 9/5/19 across
# the slide master, every slide layout, and the notes master.

$ppDateTime = 16   # ppPlaceholderDate
$newDate = "9/5/19"

function Update-DateShapes {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)

        $phType = $null
        try { $phType = $shape.PlaceholderFormat.Type } catch { $phType = $null }

        if ($phType -eq $ppDateTime) {
            if ($shape.HasTextFrame) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$p = $ppt.ActivePresentation

# Slide master.
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateShapes $layouts.Item($L).Shapes
}

# Notes master.
Update-DateShapes $p.NotesMaster.Shapes
